$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set columns C and D (Date/Time) to Text format to preserve literal strings
$ws.Range("C202:D213").NumberFormat = "@"

$ws.Range("A202").Value = "q"
$ws.Range("B202").Value = "LOG-IN"
$ws.Range("C202").Value = "02/09/24"
$ws.Range("D202").Value = "12:48:06"

$ws.Range("A203").Value = "q"
$ws.Range("B203").Value = "`nAdded item barcode no: 6454"
$ws.Range("C203").Value = "02/09/24"
$ws.Range("D203").Value = "13:39:14"

$ws.Range("A204").Value = "q"
$ws.Range("B204").Value = "Item checked out to FSL barcode no: 6454"
$ws.Range("C204").Value = "02/09/24"
$ws.Range("D204").Value = "13:40:15"

$ws.Range("A205").Value = "q"
$ws.Range("B205").Value = "Item checked in from FSL barcode no: 6454"
$ws.Range("C205").Value = "02/09/24"
$ws.Range("D205").Value = "13:40:52"

$ws.Range("A206").Value = "q"
$ws.Range("B206").Value = "LOG-OUT"
$ws.Range("C206").Value = "02/09/24"
$ws.Range("D206").Value = "13:48:20"

$ws.Range("A207").Value = "q"
$ws.Range("B207").Value = "LOG-IN"
$ws.Range("C207").Value = "02/09/24"
$ws.Range("D207").Value = "21:56:59"

$ws.Range("A208").Value = "q"
$ws.Range("B208").Value = "LOG-IN"
$ws.Range("C208").Value = "02/09/24"
$ws.Range("D208").Value = "22:01:48"

$ws.Range("A209").Value = "q"
$ws.Range("B209").Value = "LOG-IN"
$ws.Range("C209").Value = "02/09/24"
$ws.Range("D209").Value = "22:02:31"

$ws.Range("A210").Value = "q"
$ws.Range("B210").Value = "LOG-IN"
$ws.Range("C210").Value = "02/09/24"
$ws.Range("D210").Value = "22:05:16"

$ws.Range("A211").Value = "q"
$ws.Range("B211").Value = "Item checked out to FSL barcode no: 1003"
$ws.Range("C211").Value = "02/09/24"
$ws.Range("D211").Value = "22:07:20"

$ws.Range("A212").Value = "q"
$ws.Range("B212").Value = "Item checked in from FSL barcode no: 1003"
$ws.Range("C212").Value = "02/09/24"
$ws.Range("D212").Value = "22:07:45"

$ws.Range("A213").Value = "q"
$ws.Range("B213").Value = "LOG-IN"
$ws.Range("C213").Value = "02/09/24"
$ws.Range("D213").Value = "22:12:27"
